$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 246.66667
$ws.Range("I2").Value = 203.33333
$ws.Range("K2").Value = 203.33333
$ws.Range("M2").Value = -90.33332999999999
$ws.Range("H15").Value = 2761.3157
$ws.Range("I15").Value = 2761.3157
$ws.Range("K15").Value = 8283.947100000001
$ws.Range("M15").Value = -8114.947100000001
$ws.Range("H41").Value = 497.45456
$ws.Range("J41").Value = 629.4
$ws.Range("L41").Value = 629.4
$ws.Range("N41").Value = -1509.4
$ws.Range("H86").Value = 1276.875
$ws.Range("I86").Value = 1608.75
$ws.Range("J86").Value = 945
$ws.Range("K86").Value = 1608.75
$ws.Range("L86").Value = 945
$ws.Range("M86").Value = -485.75
$ws.Range("N86").Value = -3191
$ws.Range("H89").Value = 1276.875
$ws.Range("I89").Value = 1608.75
$ws.Range("J89").Value = 945
$ws.Range("K89").Value = 8043.75
$ws.Range("L89").Value = 4725
$ws.Range("M89").Value = -2427.75
$ws.Range("N89").Value = -15957
$ws.Range("H137").Value = 2403
$ws.Range("I137").Value = 2350
$ws.Range("K137").Value = 7050
$ws.Range("M137").Value = -4500
$ws.Range("H138").Value = 3525.0227
$ws.Range("I138").Value = 6999.143
$ws.Range("J138").Value = 2867.7568
$ws.Range("K138").Value = 20997.429
$ws.Range("L138").Value = 8603.270400000001
$ws.Range("M138").Value = -15857.429
$ws.Range("N138").Value = -18883.2704

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3465.836
$ws.Range("I32").Value = 2879.125
$ws.Range("J32").Value = 6452.727
$ws.Range("K32").Value = 2879.125
$ws.Range("L32").Value = 6452.727
$ws.Range("M32").Value = -2592.125
$ws.Range("N32").Value = -7026.727
$ws.Range("H45").Value = 9001313
$ws.Range("I45").Value = 10001379
$ws.Range("K45").Value = 10001379
$ws.Range("M45").Value = -10001002
$ws.Range("H61").Value = 4419.6665
$ws.Range("I61").Value = 1250
$ws.Range("J61").Value = 5325.2856
$ws.Range("K61").Value = 1250
$ws.Range("L61").Value = 5325.2856
$ws.Range("M61").Value = -1038
$ws.Range("N61").Value = -5749.2856
$ws.Range("H136").Value = 4419.6665
$ws.Range("I136").Value = 1250
$ws.Range("J136").Value = 5325.2856
$ws.Range("K136").Value = 3750
$ws.Range("L136").Value = 15975.8568
$ws.Range("M136").Value = -1200
$ws.Range("N136").Value = -21075.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 22291.2
$ws.Range("I11").Value = 824.25
$ws.Range("K11").Value = 824.25
$ws.Range("M11").Value = -684.25
$ws.Range("H107").Value = 5685.1113
$ws.Range("I107").Value = 5293
$ws.Range("J107").Value = 5998.8
$ws.Range("K107").Value = 5293
$ws.Range("L107").Value = 5998.8
$ws.Range("M107").Value = -3373
$ws.Range("N107").Value = -9838.799999999999
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 14452.75
$ws.Range("I134").Value = 14452.75
$ws.Range("K134").Value = 43358.25
$ws.Range("M134").Value = -40823.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5682893.5
$ws.Range("I22").Value = 825.8
$ws.Range("J22").Value = 10417950
$ws.Range("K22").Value = 825.8
$ws.Range("L22").Value = 10417950
$ws.Range("M22").Value = -475.8
$ws.Range("N22").Value = -10418650
$ws.Range("H31").Value = 1589.4902
$ws.Range("I31").Value = 899.4666999999999
$ws.Range("J31").Value = 1877
$ws.Range("K31").Value = 899.4666999999999
$ws.Range("L31").Value = 1877
$ws.Range("M31").Value = -604.4666999999999
$ws.Range("N31").Value = -2467
$ws.Range("H34").Value = 1589.4902
$ws.Range("I34").Value = 899.4666999999999
$ws.Range("J34").Value = 1877
$ws.Range("K34").Value = 899.4666999999999
$ws.Range("L34").Value = 1877
$ws.Range("M34").Value = -697.4666999999999
$ws.Range("N34").Value = -2281
$ws.Range("H107").Value = 551.3570999999999
$ws.Range("I107").Value = 404
$ws.Range("K107").Value = 404
$ws.Range("M107").Value = 1516

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 4000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H56").Value = 776225.25
$ws.Range("I56").Value = 776225.25
$ws.Range("K56").Value = 776225.25
$ws.Range("M56").Value = -775695.25
$ws.Range("H68").Value = 2255.8108
$ws.Range("J68").Value = 2890.1738
$ws.Range("L68").Value = 8670.5214
$ws.Range("N68").Value = -10292.5214
$ws.Range("H71").Value = 2255.8108
$ws.Range("J71").Value = 2890.1738
$ws.Range("L71").Value = 26011.5642
$ws.Range("N71").Value = -34123.5642
$ws.Range("H103").Value = 22098.584
$ws.Range("I103").Value = 31456.5
$ws.Range("J103").Value = 3382.75
$ws.Range("K103").Value = 94369.5
$ws.Range("L103").Value = 10148.25
$ws.Range("M103").Value = -93490.5
$ws.Range("N103").Value = -11906.25
$ws.Range("H107").Value = 1334.6786
$ws.Range("I107").Value = 1075
$ws.Range("J107").Value = 1377.9584
$ws.Range("K107").Value = 3225
$ws.Range("L107").Value = 4133.8752
$ws.Range("M107").Value = -1305
$ws.Range("N107").Value = -7973.8752
$ws.Range("H126").Value = 6000
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -27880
$ws.Range("H131").Value = 10885872
$ws.Range("J131").Value = 17024.227
$ws.Range("L131").Value = 51072.681
$ws.Range("N131").Value = -61152.681
$ws.Range("H132").Value = 1890.4
$ws.Range("I132").Value = 1367
$ws.Range("J132").Value = 2021.25
$ws.Range("K132").Value = 12303
$ws.Range("L132").Value = 18191.25
$ws.Range("M132").Value = -9773
$ws.Range("N132").Value = -23251.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4999.6
$ws.Range("I70").Value = 5332.6665
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 5332.6665
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -5062.6665
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 4999.6
$ws.Range("I73").Value = 5332.6665
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 5332.6665
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -4396.6665
$ws.Range("N73").Value = -6372
$ws.Range("H97").Value = 1497
$ws.Range("I97").Value = 1149.1333
$ws.Range("J97").Value = 2149.25
$ws.Range("K97").Value = 1149.1333
$ws.Range("L97").Value = 2149.25
$ws.Range("M97").Value = -653.1333
$ws.Range("N97").Value = -3141.25
$ws.Range("H132").Value = 2406627
$ws.Range("I132").Value = 3207335
$ws.Range("J132").Value = 4503.25
$ws.Range("K132").Value = 9622005
$ws.Range("L132").Value = 13509.75
$ws.Range("M132").Value = -9619475
$ws.Range("N132").Value = -18569.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5413.067
$ws.Range("J40").Value = 9799.799999999999
$ws.Range("L40").Value = 9799.799999999999
$ws.Range("N40").Value = -10071.8
$ws.Range("H46").Value = 929.2
$ws.Range("I46").Value = 447.5
$ws.Range("J46").Value = 1049.625
$ws.Range("K46").Value = 447.5
$ws.Range("L46").Value = 1049.625
$ws.Range("M46").Value = -259.5
$ws.Range("N46").Value = -1425.625
$ws.Range("H61").Value = 2523.2942
$ws.Range("I61").Value = 2074.75
$ws.Range("J61").Value = 3599.8
$ws.Range("K61").Value = 2074.75
$ws.Range("L61").Value = 3599.8
$ws.Range("M61").Value = -1872.75
$ws.Range("N61").Value = -4003.8
$ws.Range("H113").Value = 2523.2942
$ws.Range("I113").Value = 2074.75
$ws.Range("J113").Value = 3599.8
$ws.Range("K113").Value = 2074.75
$ws.Range("L113").Value = 3599.8
$ws.Range("M113").Value = 95.25
$ws.Range("N113").Value = -7939.8
$ws.Range("H136").Value = 3856.1562
$ws.Range("I136").Value = 2674.1428
$ws.Range("K136").Value = 8022.428400000001
$ws.Range("M136").Value = -5472.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 603.3
$ws.Range("I132").Value = 603.3
$ws.Range("K132").Value = 1809.9
$ws.Range("M132").Value = 720.1000000000001
